# Natmi following Dr Hou advice
# Replace/extend the LR-pairs data (rows 2-13, columns A-T) on the active sheet
# with the re-computed ligand/receptor statistics (sending cluster x target
# cluster combinations for FAPs and sCs against the full set of target
# clusters: ECs, FAPs, M1, M2, Neutro, sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  ,@("FAPs","Has2","Hmmr","ECs",3,1,48.25514733333333,144.765442,0.9340796096783166,0.9340796096783167,2,0.6666666666666666,1.375678666666666,4.127035999999999,0.03616159583349148,0.03616159583349148,66.38357674332354,597.452190689912,0.03377780932149276,0.03377780932149277)
  ,@("FAPs","Has2","Hmmr","FAPs",3,1,48.25514733333333,144.765442,0.9340796096783166,0.9340796096783167,3,1,25.73962066666667,77.218862,0.6766011438635753,0.6766011438635753,1242.069187574112,11178.622688167,0.6319993323679909,0.631999332367991)
  ,@("FAPs","Has2","Hmmr","M1",3,1,48.25514733333333,144.765442,0.9340796096783166,0.9340796096783167,2,0.6666666666666666,2.416767333333333,7.250302,0.06352803576095653,0.06352803576095653,116.6214637403871,1049.593173663484,0.05934024284722442,0.05934024284722442)
  ,@("FAPs","Has2","Hmmr","M2",3,1,48.25514733333333,144.765442,0.9340796096783166,0.9340796096783167,3,1,1.172619333333333,3.517858,0.03082390344925866,0.03082390344925867,56.58491869591511,509.264268263236,0.02879197970264565,0.02879197970264566)
  ,@("FAPs","Has2","Hmmr","Neutro",3,1,48.25514733333333,144.765442,0.9340796096783166,0.9340796096783167,3,1,1.025619333333333,3.076858,0.02695980733704406,0.02695980733704405,49.49141203791512,445.422708341236,0.02518260631438873,0.02518260631438873)
  ,@("FAPs","Has2","Hmmr","sCs",3,1,48.25514733333333,144.765442,0.9340796096783166,0.9340796096783167,3,1,6.312226666666667,18.93668,0.1659255137556739,0.1659255137556739,304.5974278013956,2741.37685021256,0.1549876391245741,0.1549876391245741)
  ,@("sCs","Has2","Hmmr","ECs",3,1,3.405489333333334,10.216468,0.06592039032168336,0.06592039032168336,2,0.6666666666666666,1.375678666666666,4.127035999999999,0.03616159583349148,0.03616159583349148,4.684859025427556,42.163731228848,0.002383786511998717,0.002383786511998717)
  ,@("sCs","Has2","Hmmr","FAPs",3,1,3.405489333333334,10.216468,0.06592039032168336,0.06592039032168336,3,1,25.73962066666667,77.218862,0.6766011438635753,0.6766011438635753,87.65600362437957,788.9040326194161,0.04460181149558431,0.04460181149558431)
  ,@("sCs","Has2","Hmmr","M1",3,1,3.405489333333334,10.216468,0.06592039032168336,0.06592039032168336,2,0.6666666666666666,2.416767333333333,7.250302,0.06352803576095653,0.06352803576095653,8.230275374815111,74.072478373336,0.004187792913732113,0.004187792913732113)
  ,@("sCs","Has2","Hmmr","M2",3,1,3.405489333333334,10.216468,0.06592039032168336,0.06592039032168336,3,1,1.172619333333333,3.517858,0.03082390344925866,0.03082390344925867,3.993342631727111,35.940083685544,0.002031923746613013,0.002031923746613013)
  ,@("sCs","Has2","Hmmr","Neutro",3,1,3.405489333333334,10.216468,0.06592039032168336,0.06592039032168336,3,1,1.025619333333333,3.076858,0.02695980733704406,0.02695980733704405,3.492735699727112,31.434621297544,0.001777201022655327,0.001777201022655327)
  ,@("sCs","Has2","Hmmr","sCs",3,1,3.405489333333334,10.216468,0.06592039032168336,0.06592039032168336,3,1,6.312226666666667,18.93668,0.1659255137556739,0.1659255137556739,21.49622058291556,193.46598524624,0.01093787463109987,0.01093787463109987)
)

$nRows = $rows.Count
$nCols = $rows[0].Count
$data = New-Object 'object[,]' $nRows,$nCols
for ($r = 0; $r -lt $nRows; $r++) {
  for ($c = 0; $c -lt $nCols; $c++) {
    $data[$r,$c] = $rows[$r][$c]
  }
}

$ws.Range("A2:T13").Value = $data
